$d = $word.ActiveDocument

# Paragraph 1 is the Title paragraph: "Module 1" -> "Wireshark Activity"
$titlePara = $d.Paragraphs(1).Range
$titlePara.Find.Execute("Module 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wireshark Activity", 2) | Out-Null

# Paragraph 2 is the Subtitle paragraph: "Wireshark Activity" -> "Wagner Module 1"
$subtitlePara = $d.Paragraphs(2).Range
$subtitlePara.Find.Execute("Wireshark Activity", $true, $false, $false, $false, $false,
                            $true, 1, $false, "Wagner Module 1", 2) | Out-Null
